# Adapt column header formatting to respective input file names.
#   *_old  -> *_FV2404
#   *_new  -> *_FV2410
# Freeze the header row and wrap the data range in an Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- 1. Rename the header row (row 1) -------------------------------------
$oldSuffix = "_old"
$newSuffixFor2404 = "_FV2404"
$newSuffix = "_new"
$newSuffixFor2410 = "_FV2410"

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -eq $null) { continue }
    $val = [string]$val
    if ($val.EndsWith($oldSuffix)) {
        $base = $val.Substring(0, $val.Length - $oldSuffix.Length)
        $cell.Value = $base + $newSuffixFor2404
    } elseif ($val.EndsWith($newSuffix)) {
        $base = $val.Substring(0, $val.Length - $newSuffix.Length)
        $cell.Value = $base + $newSuffixFor2410
    }
}

# --- 2. Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the data range into an Excel Table -----------------------------
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Output "done"
